$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (student 424346 scan) is being replaced by what used to be row 3's
# scan (student 676767, logged 2 seconds later). Copy the cells instead of
# assigning literal values so the text-stored-as-number formatting of the
# Student ID / Log Time columns is preserved exactly as Excel had it.
$ws.Range("A3").Copy($ws.Range("A2"))
$ws.Range("D3").Copy($ws.Range("D2"))

# The old row 3 (now duplicated into row 2) is removed, shifting everything
# below it up and shrinking the used range from A1:F3 to A1:F2.
$ws.Rows("3:3").Delete()
